$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $searchText, $replaceText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($searchText)
    if ($idx -lt 0) {
        throw "Text not found in paragraph $paraIndex (was: $full)"
    }
    $start = $p.Range.Start + $idx
    $end = $start + $searchText.Length
    $rr = $d.Range($start, $end)
    $rr.Text = $replaceText
}

function Set-ParagraphText($paraIndex, $newText) {
    # Always force a genuine content-length change so the engine actually
    # rewrites/collapses the paragraph's runs (no-op same-text assignments
    # are ignored and leave old run/proofErr structure untouched).
    $p = $d.Paragraphs.Item($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    $rr = $d.Range($start, $end)
    $rr.Text = $newText + "ZZZTMPZZZ"

    $p2 = $d.Paragraphs.Item($paraIndex)
    $full2 = $p2.Range.Text
    $markerIdx = $full2.IndexOf("ZZZTMPZZZ")
    $markerStart = $p2.Range.Start + $markerIdx
    $markerEnd = $markerStart + 9
    $rr2 = $d.Range($markerStart, $markerEnd)
    $rr2.Delete()
}

$ppoClause = "ppo_type == “nondomestic” and respondent_is_minor and (not respondent_is_emancipated_minor)"
$trimmedClause = "respondent_is_minor and (not respondent_is_emancipated_minor)"

# 1) "Other Pending Actions" continuation condition (paragraph 4 / B2a):
#    drop the leading ppo_type clause, keep "{% if " and the " %}B2a{% else %}C3a{% endif %}" tail intact.
Replace-InParagraph 4 $ppoClause $trimmedClause

# 2) "Name of judge" block for `action` (paragraphs 23-24): collapse the
#    proofErr-wrapped runs back into single clean runs.
Set-ParagraphText 23 "{%p if action.judge %}"
Set-ParagraphText 24 "Name of judge: {{ action.judge }}"

# 3) "Orders and Judgments" continuation condition (paragraph 32 / B2b):
#    same ppo_type clause removal as step 1.
Replace-InParagraph 32 $ppoClause $trimmedClause

# 4) "Relief Requested - Other Relief" continuation condition (paragraph 60 / D4e):
#    here the clause removal also swallows the "if " keyword, per the target diff.
Replace-InParagraph 60 "if $ppoClause " "$trimmedClause "
